$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (1-11) with new type/name/label values ---
$ws.Range("A1").Value = "type"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "label"

$ws.Range("A2").Value = "start"
$ws.Range("B2").Value = "start"

$ws.Range("A3").Value = "end"
$ws.Range("B3").Value = "end"

$ws.Range("A4").Value = "note"
$ws.Range("B4").Value = "form_name"
$ws.Range("C4").Value = "Form Name"

$ws.Range("A5").Value = "note"
$ws.Range("B5").Value = "form_version"
$ws.Range("C5").Value = "Version"

$ws.Range("A6").Value = "text"
$ws.Range("B6").Value = "code"
$ws.Range("C6").Value = "What is the 3 letter site Code?"

$ws.Range("A7").Value = "geopoint"
$ws.Range("B7").Value = "gps_code"
$ws.Range("C7").Value = "Select a single GPS point to represent this 3 letter site code"

$ws.Range("A8").Value = "begin_group"
$ws.Range("B8").Value = "group_000"
$ws.Range("C8").Value = "indicates a barcode and data is coming"

$ws.Range("A9").Value = "barcode"
$ws.Range("B9").Value = "barcode_bag_000"
$ws.Range("C9").Value = "Barcode for decomp bag"

$ws.Range("A10").Value = "decimal"
$ws.Range("B10").Value = "fresh_bag_wt_grams_000"
$ws.Range("C10").Value = "Total Fresh Wt (grams)"

$ws.Range("A11").Value = "decimal"
$ws.Range("B11").Value = "pre_bag_wt_grams_000"
$ws.Range("C11").Value = "Decomp Bag Pre Wt"

# --- Row 12: only B12 remains (old A12/C12 begin_group data is gone) ---
$ws.Range("A12").ClearContents()
$ws.Range("B12").Value = "dry_bag_wt_grams_000"
$ws.Range("C12").ClearContents()

# --- Row 13 is fully removed ---
$ws.Range("A13").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()

# --- Row 14 now holds a note instead of the old decimal/dry_wt_001 row ---
$ws.Range("A14").Value = "Barcodes coming in always start with 'barcode_'"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# --- New rows 15-18 ---
$ws.Range("A15").Value = "Site/farm/code"
$ws.Range("B15").Value = "code"

$ws.Range("A16").Value = "Groups:"
$ws.Range("B16").Value = "group_000"

$ws.Range("A17").Value = "Yield:"
$ws.Range("B17").Value = "barcode_yield_000"

$ws.Range("B18").Value = "yield_wt_000"

# --- Highlight fills (yellow), matches style used elsewhere in the sheet ---
$ws.Range("B12").Interior.Color = 65535
$ws.Range("B18").Interior.Color = 65535

# --- Column widths (best effort / bestfit-style; column C keeps its original width) ---
$ws.Columns("A").ColumnWidth = 40.33
$ws.Columns("B").ColumnWidth = 21.65

# --- Selection moves to B14 ---
$ws.Range("B14").Select() | Out-Null
